{"js": "// Update the division-fact worksheet: each \"AA\u00f7B=\" cell in the table is\n// replaced with a new \"CC\u00f7D=\" problem (commit: \"Update master to output\n// generated at 9a8706d\"). Several old problems repeat (e.g. \"83\u00f73=\" and\n// \"45\u00f75=\" each appear twice), so matches are consumed strictly in\n// document order and mapped 1:1 onto the corresponding replacement list.\nconst replacements = [\n  [\"10\u00f78=\", [\"81\u00f78=\"]],\n  [\"61\u00f76=\", [\"57\u00f74=\"]],\n  [\"27\u00f78=\", [\"98\u00f75=\"]],\n  [\"48\u00f74=\", [\"23\u00f74=\"]],\n  [\"43\u00f77=\", [\"12\u00f75=\"]],\n  [\"83\u00f73=\", [\"80\u00f72=\", \"61\u00f75=\"]],\n  [\"30\u00f74=\", [\"40\u00f74=\"]],\n  [\"87\u00f72=\", [\"39\u00f77=\"]],\n  [\"45\u00f75=\", [\"14\u00f77=\", \"69\u00f77=\"]],\n  [\"58\u00f78=\", [\"67\u00f72=\"]],\n  [\"88\u00f73=\", [\"41\u00f76=\"]],\n  [\"32\u00f75=\", [\"40\u00f75=\"]],\n  [\"69\u00f74=\", [\"54\u00f75=\"]],\n  [\"94\u00f73=\", [\"97\u00f74=\"]],\n  [\"29\u00f79=\", [\"52\u00f74=\"]],\n  [\"95\u00f73=\", [\"53\u00f78=\"]],\n  [\"97\u00f75=\", [\"50\u00f77=\"]],\n  [\"65\u00f76=\", [\"76\u00f73=\"]],\n  [\"10\u00f76=\", [\"89\u00f72=\"]],\n  [\"10\u00f72=\", [\"67\u00f72=\"]],\n  [\"44\u00f79=\", [\"61\u00f75=\"]],\n  [\"57\u00f72=\", [\"63\u00f73=\"]],\n  [\"21\u00f79=\", [\"65\u00f74=\"]],\n];\n\nfor (const [oldText, newTexts] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length !== newTexts.length) {\n    throw new Error(\n      `expected ${newTexts.length} match(es) for \"${oldText}\", found ${found.items.length}`\n    );\n  }\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newTexts[i], \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-fact worksheet: each \"AA\u00f7B=\" cell in the table is\n# replaced with a new \"CC\u00f7D=\" problem (commit: \"Update master to output\n# generated at 9a8706d\"). The table is a fixed 20-row x 5-col grid where\n# only every 4th row (1,5,9,13,17) carries data, so cells are addressed\n# directly by (row, col) -- this sidesteps the fact that some old values\n# (\"83\u00f73=\", \"45\u00f75=\") repeat and need different replacements per occurrence.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n  @{ Row = 1;  Col = 1; Old = \"10\u00f78=\"; New = \"81\u00f78=\" },\n  @{ Row = 1;  Col = 2; Old = \"61\u00f76=\"; New = \"57\u00f74=\" },\n  @{ Row = 1;  Col = 3; Old = \"27\u00f78=\"; New = \"98\u00f75=\" },\n  @{ Row = 1;  Col = 4; Old = \"48\u00f74=\"; New = \"23\u00f74=\" },\n  @{ Row = 1;  Col = 5; Old = \"43\u00f77=\"; New = \"12\u00f75=\" },\n  @{ Row = 5;  Col = 1; Old = \"83\u00f73=\"; New = \"80\u00f72=\" },\n  @{ Row = 5;  Col = 2; Old = \"30\u00f74=\"; New = \"40\u00f74=\" },\n  @{ Row = 5;  Col = 3; Old = \"87\u00f72=\"; New = \"39\u00f77=\" },\n  @{ Row = 5;  Col = 4; Old = \"45\u00f75=\"; New = \"14\u00f77=\" },\n  @{ Row = 5;  Col = 5; Old = \"58\u00f78=\"; New = \"67\u00f72=\" },\n  @{ Row = 9;  Col = 1; Old = \"88\u00f73=\"; New = \"41\u00f76=\" },\n  @{ Row = 9;  Col = 2; Old = \"83\u00f73=\"; New = \"61\u00f75=\" },\n  @{ Row = 9;  Col = 3; Old = \"32\u00f75=\"; New = \"40\u00f75=\" },\n  @{ Row = 9;  Col = 4; Old = \"69\u00f74=\"; New = \"54\u00f75=\" },\n  @{ Row = 9;  Col = 5; Old = \"94\u00f73=\"; New = \"97\u00f74=\" },\n  @{ Row = 13; Col = 1; Old = \"29\u00f79=\"; New = \"52\u00f74=\" },\n  @{ Row = 13; Col = 2; Old = \"95\u00f73=\"; New = \"53\u00f78=\" },\n  @{ Row = 13; Col = 3; Old = \"97\u00f75=\"; New = \"50\u00f77=\" },\n  @{ Row = 13; Col = 4; Old = \"65\u00f76=\"; New = \"76\u00f73=\" },\n  @{ Row = 13; Col = 5; Old = \"10\u00f76=\"; New = \"89\u00f72=\" },\n  @{ Row = 17; Col = 1; Old = \"10\u00f72=\"; New = \"67\u00f72=\" },\n  @{ Row = 17; Col = 2; Old = \"44\u00f79=\"; New = \"61\u00f75=\" },\n  @{ Row = 17; Col = 3; Old = \"57\u00f72=\"; New = \"63\u00f73=\" },\n  @{ Row = 17; Col = 4; Old = \"21\u00f79=\"; New = \"65\u00f74=\" },\n  @{ Row = 17; Col = 5; Old = \"45\u00f75=\"; New = \"69\u00f77=\" }\n)\n\nforeach ($edit in $edits) {\n  $cell = $t.Cell($edit.Row, $edit.Col)\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $edit.Old) {\n    throw \"cell ($($edit.Row),$($edit.Col)): expected '$($edit.Old)' but found '$current'\"\n  }\n  $cell.Range.Text = $edit.New\n}\n\nWrite-Output \"done\"\n"}
